$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Timings" column header is expanded into a fuller example showing the
# expected format, e.g. "Timings  Like (1:00 PM To 3:00 PM)".
$ws.Range("B1").Value = "Timings  Like (1:00 PM To 3:00 PM)"

# Column B is widened so the longer header text fits/reads well
# (target stored width ~21.66 characters).
$ws.Columns("B").ColumnWidth = 20.8333333

# The saved cursor/selection moves from E10 to E8.
$ws.Range("E8").Select() | Out-Null
